$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B: all zero for rows 2-9
$ws.Range("B2:B9").Value = 0

# Column D: all zero for rows 2-9
$ws.Range("D2:D9").Value = 0

# Column C: zero by default, with specific overrides
$ws.Range("C2:C9").Value = 0
$ws.Range("C2").Value = 0.7844137198105535
$ws.Range("C4").Value = -0.7404339077487648
$ws.Range("C6").Value = -0.6919333072105328
$ws.Range("C8").Value = 0.7757279151612556
